# geração de análises seriais
# Swap tied-ranking category labels on two sheets.

$wb = $excel.ActiveWorkbook

# Sheet "max-arrecad": rows 9 and 10 are tied (415470.057) -> swap order of labels
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A9").Value = "fiq"
$wsMax.Range("A10").Value = "folclore"

# Sheet "tx-sucesso": rows 15 and 16 are tied (66.2) -> swap order of labels
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A15").Value = "politica"
$wsTx.Range("A16").Value = "erotismo"
